$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header labels - swap A1/B1 and D1/E1
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "bedrooms_2"

# Row 2 - swap D2/E2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

# Row 3 - move the 1 from B3 to D3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1

# Row 4 - swap A4/B4
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 1

# Row 6 - move the 1 from E6 to A6
$ws.Range("A6").Value = 1
$ws.Range("E6").Value = 0
